# Adding company cidr to template
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: DC_vlans_optional - insert 11 new rows of CIDR data after
# row 6, then append 8 more rows of CIDR data at the bottom.
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("DC_vlans_optional")

# Insert 11 blank rows starting at row 7 (shifts old rows 7-37 down to 18-48)
$ws4.Rows.Item(7).Resize(11).Insert()

# Fill newly inserted rows 7-17 with new CIDR entries
$ws4.Range("A7").Value = "172.28.40.0/29"
$ws4.Range("B7").Value = "Wireless_net1"
$ws4.Range("A8").Value = "172.28.40.8/29"
$ws4.Range("B8").Value = "DNAC_mgmt"
$ws4.Range("A9").Value = "172.28.40.16/29"
$ws4.Range("B9").Value = "Wireless_net3"
$ws4.Range("A10").Value = "172.28.40.24/29"
$ws4.Range("B10").Value = "SilverPeak_mgmt"
$ws4.Range("A11").Value = "172.28.44.0/24"
$ws4.Range("B11").Value = "EN_Wireless_Users_44"
$ws4.Range("A12").Value = "10.248.0.0/16"
$ws4.Range("B12").Value = "SDA_POV"
$ws4.Range("A13").Value = "10.123.0.0/16"
$ws4.Range("B13").Value = "SDA_POV"
$ws4.Range("A14").Value = "10.124.0.0/16"
$ws4.Range("B14").Value = "SDA_POV"
$ws4.Range("A15").Value = "10.125.0.0/16"
$ws4.Range("B15").Value = "SDA_POV"
$ws4.Range("A16").Value = "10.126.0.0/16"
$ws4.Range("B16").Value = "SDA_POV"
$ws4.Range("A17").Value = "10.127.0.0/16"
$ws4.Range("B17").Value = "SDA_POV"

# Append 8 new rows (49-56) after the existing data (now ending at row 48)
$ws4.Range("A49").Value = "172.17.25.0/24"
$ws4.Range("B49").Value = "EMEA_servers"
$ws4.Range("A50").Value = "172.18.153.0/24"
$ws4.Range("B50").Value = "Singapore_servers"
$ws4.Range("A51").Value = "172.19.120.0/24"
$ws4.Range("B51").Value = "APAC_servers"
$ws4.Range("A52").Value = "172.23.107.0/24"
$ws4.Range("B52").Value = "APAC_servers"
$ws4.Range("A53").Value = "172.23.107.0/24"
$ws4.Range("B53").Value = "APAC_servers"
$ws4.Range("A54").Value = "172.24.32.96/29"
$ws4.Range("B54").Value = "ATL_servers"
$ws4.Range("A55").Value = "172.30.254.0/24"
$ws4.Range("B55").Value = "APIC-EM_System"
$ws4.Range("A56").Value = "172.31.150.0/23"
$ws4.Range("B56").Value = "EMEA_servers"

# Update the sheet's selection to B7 (matches the saved view state)
[void]$ws4.Range("B7").Select()

# ---------------------------------------------------------------
# Sheet: Company_cidr - cosmetic view/column width refresh
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Company_cidr")
$ws1.Columns.Item(1).ColumnWidth = 23.3

# Re-activate the first sheet so it remains the active/selected tab
# (selecting on the other sheet above would otherwise move the
# "active tab" flag away from it).
[void]$ws1.Activate()
